$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.990.68'
$ws.Range('E2').Value = '  +0.93%  '
$ws.Range('D3').Value = '1.556.52'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.34'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.68'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.248'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('E10').Value = '  +1.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0862'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = '1.779.71'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = '1.558.12'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.72'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +1.49%  '
$ws.Range('E15').Value = '  +1.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.00'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +1.41%  '
$ws.Range('D17').Value = '26.999.70'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '216.04'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.27'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.19'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +2.96%  '
$ws.Range('E24').Value = '  -0.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.41'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('E26').Value = '  +2.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.92'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('E29').Value = '  +1.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0463'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('E32').Value = '  +1.22%  '
$ws.Range('D33').Value = '1.406.08'
$ws.Range('E33').Value = '  +4.55%  '
$ws.Range('E34').Value = '  +3.15%  '
$ws.Range('E35').Value = '  +3.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.958'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +3.42%  '
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('E38').Value = '  +1.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.523'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.810'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +1.24%  '
$ws.Range('E41').Value = '  +0.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.989'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('E43').Value = '  +3.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.50'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -4.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.97'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +1.91%  '
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').Value = '1.693.41'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.29'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0957'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('E51').Value = '  +0.53%  '
